$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.707.32'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '1.732.42'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("D4").Value = '''0.9982'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '''242.31'
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").Value = '''0.4931'
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("D8").Value = '''0.2624'
$ws.Range("E8").Value = '  +0.42%  '
$ws.Range("D9").Value = '''0.06218'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").Value = '1.726.90'
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").Value = '''15.91'
$ws.Range("E11").Value = '  +3.38%  '
$ws.Range("D12").Value = '''0.06984'
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '''0.6104'
$ws.Range("E13").Value = '  +2.18%  '
$ws.Range("D14").Value = '''4.502'
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").Value = '''0.9984'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = '26.514.75'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").Value = '''0.9986'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = '''0.000007201'
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '''11.41'
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("D21").Value = '1.949.35'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = '''4.470'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '''8.554'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '''5.090'
$ws.Range("E24").Value = '  -1.54%  '
$ws.Range("D25").Value = '''138.50'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("D26").Value = '''15.35'
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("D27").Value = '''1.768'
$ws.Range("E27").Value = '  +3.36%  '
$ws.Range("E28").Value = '  -2.05%  '
$ws.Range("D29").Value = '''106.67'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("D30").Value = '''3.934'
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("D31").Value = '''0.07982'
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("D32").Value = '''3.664'
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("D35").Value = '''1.002'
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("D36").Value = '''0.6228'
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("D37").Value = '''0.9419'
$ws.Range("E37").Value = '  +4.03%  '
$ws.Range("D38").Value = '''2.041'
$ws.Range("E38").Value = '  +2.91%  '
$ws.Range("D39").Value = '''2.420'
$ws.Range("E39").Value = '  +0.66%  '
$ws.Range("E41").Value = '  +1.70%  '
$ws.Range("D42").Value = '''5.571'
$ws.Range("E42").Value = '  +3.19%  '
$ws.Range("D43").Value = '''99.46'
$ws.Range("D44").Value = '''0.3856'
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("D45").Value = '''6.923'
$ws.Range("E45").Value = '  +2.99%  '
$ws.Range("E46").Value = '  +1.05%  '
$ws.Range("D47").Value = '''0.05379'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").Value = '''7.898'
$ws.Range("E48").Value = '  +2.52%  '
$ws.Range("D49").Value = '''30.25'
$ws.Range("E49").Value = '  +0.40%  '
$ws.Range("D50").Value = '''51.66'
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("D51").Value = '''1.229'
$ws.Range("E51").Value = '  -1.06%  '
